$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8181886076927185
$ws.Range("B1").Value = 1.147456288337708
$ws.Range("C1").Value = 1.597726225852966
$ws.Range("D1").Value = 4.71535062789917
$ws.Range("E1").Value = 2.240630388259888
